# Applies the diff: updated Nmap scan timestamps/latencies and reordered
# / inserted Metasploit exploit log lines (with updated session port/time).

$d = $word.ActiveDocument

# --- Nmap scan block -------------------------------------------------

# 1. "# Nmap 7.80 scan initiated ..." timestamp
$d.Content.Find.Execute(
    "# Nmap 7.80 scan initiated Sun Jul  7 03:31:54 2024 as:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "# Nmap 7.80 scan initiated Sun Jul  7 03:38:31 2024 as:", 2) | Out-Null

# 2. Latency for 10.33.102.225 (first "Host is up" occurrence)
$d.Content.Find.Execute(
    "Host is up (0.00061s latency).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Host is up (0.00070s latency).", 2) | Out-Null

# 3. Latency for 10.33.102.226 (second "Host is up" occurrence)
$d.Content.Find.Execute(
    "Host is up (0.00066s latency).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Host is up (0.00076s latency).", 2) | Out-Null

# 4. "# Nmap done at ..." timestamp + scan duration
$d.Content.Find.Execute(
    "# Nmap done at Sun Jul  7 03:32:02 2024 -- 2 IP addresses (2 hosts up) scanned in 7.78 seconds",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "# Nmap done at Sun Jul  7 03:38:38 2024 -- 2 IP addresses (2 hosts up) scanned in 7.82 seconds", 2) | Out-Null

# --- Metasploit exploit log block ------------------------------------
# The original block (lines joined by manual line breaks, i.e. Chr(11)
# in Range.Text) is reordered and gains one new line ("Sending stage ..."
# right after the automatic-check line); the Meterpreter session line's
# port/timestamp and the "Sending stage"/"Command Stager" ordering also
# change.
#
# NOTE: the Find/Replace "insert" path re-runs AutoFormat's smart-quote
# substitution on the *replacement* text, which would corrupt the
# straight double quotes in the untouched
# `("set AutoCheck false" to disable)` line. Anchor on `disable)`
# (just after the quoted text) instead of the literal quote characters
# so no `"` ever appears inside a replacement string.

$vt = [char]11

$oldBlock = @(
    'disable)',
    '[+] The target appears to be vulnerable. The target is Cacti version 1.2.22',
    '[*] Trying to bruteforce an exploitable host_id and local_data_id by trying up to 500 combinations',
    '[*] Enumerating local_data_id values for host_id 1',
    '[*] Sending stage (1017704 bytes) to 10.33.102.225',
    '[*] Meterpreter session 1 opened (10.33.102.224:4444 -> 10.33.102.225:51034) at 2024-07-07 03:35:08 +0700',
    '[+] Found exploitable local_data_id 15 for host_id 1',
    '[*] Sending stage (1017704 bytes) to 10.33.102.225',
    '[*] Command Stager progress - 100.00% done (1118/1118 bytes)'
) -join $vt

$newBlock = @(
    'disable)',
    '[*] Sending stage (1017704 bytes) to 10.33.102.225',
    '[+] The target appears to be vulnerable. The target is Cacti version 1.2.22',
    '[*] Trying to bruteforce an exploitable host_id and local_data_id by trying up to 500 combinations',
    '[*] Enumerating local_data_id values for host_id 1',
    '[*] Meterpreter session 1 opened (10.33.102.224:4444 -> 10.33.102.225:59180) at 2024-07-07 03:41:44 +0700',
    '[+] Found exploitable local_data_id 15 for host_id 1',
    '[*] Command Stager progress - 100.00% done (1118/1118 bytes)',
    '[*] Sending stage (1017704 bytes) to 10.33.102.225'
) -join $vt

$d.Content.Find.Execute($oldBlock, $true, $false, $false, $false, $false, $true, 1, $false, $newBlock, 2) | Out-Null

# --- Active-sessions table: updated source port for 10.33.102.225 ----
$d.Content.Find.Execute(
    "10.33.102.224:4444 -> 10.33.102.225:51034 (172.24.0.3)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "10.33.102.224:4444 -> 10.33.102.225:59180 (172.24.0.3)", 2) | Out-Null
